$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph (search by an
# ASCII-safe substring to sidestep any accented-character matching issues).
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Docente(s)")) {
        $targetIndex = $i
        break
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

$target = $d.Paragraphs($targetIndex)

# Split off a new, empty paragraph right after the heading.
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($targetIndex + 1)
$newRange = $newPara.Range

# Fill the new paragraph via raw OOXML so the run layout matches exactly:
# one run holding the first name followed by a manual line break, and a
# second run holding the second name, all under the "ListBullet" style.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t><w:br/></w:r><w:r><w:t>5840521 - Rosa Ana Conte</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($xml)
